$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new column K "Lineage Key" ---
$ws.Range("K1").Value = "Lineage Key"
$ws.Range("K1").Font.Bold = $true

# --- Column widths for new columns K (11) and L (12) ---
$ws.Columns.Item(11).ColumnWidth = 10.2
$ws.Columns.Item(12).ColumnWidth = 18.75

# --- Data rows 2-21: add Lineage Key value (2) in column K ---
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 11).Value = 2
}

# --- Column L formulas: rebuild the INSERT statement formula including Lineage Key ---
# Set the shared block L3:L21 first so it forms its own shared-formula group.
$ws.Range("L3:L21").Formula = '="INSERT into [dbo].[Dim_Branch] ([_Source Key], [Branch Name], [Branch PINCODE], [Number of Employees], [Number of Customers], [Branch Type], [City], [State], [Country], [Lineage Key]) VALUES ("&B3&",''"&C3&"'',''"&D3&"'',"&E3&",''"&F3&"'',''"&G3&"'',''"&H3&"'',''"&I3&"'',''"&J3&"'',"&K3&");"'

# L2 is set afterwards as its own standalone formula (it no longer belongs to the shared group).
$ws.Range("L2").Formula = '="INSERT into [dbo].[Dim_Branch] ([_Source Key], [Branch Name], [Branch PINCODE], [Number of Employees], [Number of Customers], [Branch Type], [City], [State], [Country], [Lineage Key]) VALUES ("&B2&",''"&C2&"'',''"&D2&"'',"&E2&",''"&F2&"'',''"&G2&"'',''"&H2&"'',''"&I2&"'',''"&J2&"'',"&K2&");"'

# --- Sheet view: zoom + selection changes ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("K5").Select()
